$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data to append (row 53)
$row = 53

# Column A holds a date-like string ("2025-10-02"). Excel's COM layer will
# auto-convert a bare string that looks like a date into a date serial
# number when assigned to a General-formatted cell. To preserve the
# original workbook's convention of storing these values as plain text
# (inline/shared strings, not dates), force the cell to Text format before
# assigning the value, then clear the formatting override afterwards so no
# stray style index is left on the cell (matching the rest of the sheet,
# which carries no explicit cell styles).
$ws.Range("A" + $row).NumberFormat = "@"
$ws.Range("A" + $row).Value = "2025-10-02"
$ws.Range("A" + $row).ClearFormats()

$ws.Range("B" + $row).Value = "15:23:32"
$ws.Range("C" + $row).Value = "1.00 EUR = 1,774.7030"
